# Update Active_Outages.xlsx - 6/18/2025, 10:13:25 AM
#
# Refreshes the "Elapsed Duration(Hrs)" figures on several regional sheets
# and appends a newly opened outage (JED0123 / R4) to sheet R1.

$wb = $excel.ActiveWorkbook

# --- Sheet R1 -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3923:27:31"
$ws1.Range("G3").Value = "63:00:09"

$ws1.Cells.Item(6, 1).Value  = ""
$ws1.Cells.Item(6, 2).Value  = "R4"
$ws1.Cells.Item(6, 3).Value  = ""
$ws1.Cells.Item(6, 4).Value  = "JED0123"
$ws1.Cells.Item(6, 5).Value  = ""
$ws1.Cells.Item(6, 6).Value  = ""
$ws1.Cells.Item(6, 7).Value  = ""
$ws1.Cells.Item(6, 8).Value  = ""
$ws1.Cells.Item(6, 9).Value  = "SCECO"
$ws1.Cells.Item(6, 10).Value = "In progress"
$ws1.Cells.Item(6, 11).Value = ""
$ws1.Cells.Item(6, 12).Value = "Latis"

# --- Sheet R2 -------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12104:51:11"
$ws2.Range("G3").Value = "3234:34:40"
$ws2.Range("G4").Value = "472:46:14"

# --- Sheet R4 -------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2950:41:00"
$ws4.Range("G3").Value = "177:53:15"

# --- Sheet R5 -------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "424:39:59"

# --- Sheet R6 -------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "65:12:17"
